$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: mint a fresh numbered-list definition (numId=2) by applying
# default list numbering to a scratch paragraph appended at the very
# end of the document. This creates the abstractNum/num entries in
# numbering.xml that our real bullet paragraphs will reference below.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$mintPara = $d.Paragraphs.Last
$mintPara.Style = "List Paragraph"
$mintPara.Range.Text = "x"
$mintPara.Range.ListFormat.ApplyNumberDefault()

# ------------------------------------------------------------------
# Step 2: replace the whole "Programming setup" body paragraph (which
# currently holds the single-run "Python was used..." sentence and the
# _GoBack bookmark) together with the scratch paragraph from step 1
# with the full new content: the expanded intro sentence (with a bold
# "pip"), the three package bullet points (numId=2), the "Project
# structure" heading, and a trailing empty paragraph that keeps the
# _GoBack bookmark.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(11)
$scratchPara = $d.Paragraphs.Item(12)
$target = $d.Range($introPara.Range.Start, $scratchPara.Range.End)

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Python </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">(version 3) </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>was used as a programming language for this assignment. It was chosen for it’s ease of use and comprehensive list of available packages which could simplify the development of the program.</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Package-management system </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>pip</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> allowed us to use the following packages:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a6"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>NumPy – library which supports multi-dimensional arrays and matrices and also provides the functions to operate on them.</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> We also used it as a tool to import data from TXT and CSV files.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a6"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Scikit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>-learn</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>machine learning library</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> which gave us the ability to </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">calculate the </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Normalized Mutual Information</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> score.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a6"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Matplotlib</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – plotting library which was used by us </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>to visualize the clustering result.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Project structure</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$target.InsertXML($xml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
foreach ($para in $d.Paragraphs) {
    Write-Output ("[" + $para.Range.Text + "]")
}
